$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (date serial 44511 = 2021-11-11) needs to be
# inserted as the first entry of this block (row 96), pushing the existing
# rows 96-119 down to 97-120.
$ws.Rows.Item(96).Insert()

$newRow = 96
$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 44511
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100101
$ws.Cells.Item($newRow, 8).Value = "Berries"
$ws.Cells.Item($newRow, 9).Value = 100112025
$ws.Cells.Item($newRow, 10).Value = "Frutilla"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 400
$ws.Cells.Item($newRow, 14).Value = 9000
$ws.Cells.Item($newRow, 15).Value = 10000
$ws.Cells.Item($newRow, 16).Value = 9500
$ws.Cells.Item($newRow, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item($newRow, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item($newRow, 19).Value = 1357
$ws.Cells.Item($newRow, 20).Value = 7
